$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.735.32'
$ws.Range("E2").Value = '  +5.75%  '
$ws.Range("D3").Value = '2.259.19'
$ws.Range("E3").Value = '  +4.40%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.19'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("E6").Value = '  +3.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.39'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +4.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '60.19'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.106'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("D13").Value = '2.593.69'
$ws.Range("E13").Value = '  +4.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.79'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.90%  '
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '2.249.83'
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("D19").Value = '41.607.44'
$ws.Range("E19").Value = '  +5.41%  '
$ws.Range("D20").Value = '0.0₃0943'
$ws.Range("E20").Value = '  +10.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.50'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.02'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +10.17%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.148'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.65'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.56'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.38%  '
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("E32").Value = '  +7.96%  '
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.08'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.80'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0642'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.93'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.86'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.50'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("E40").Value = '  +58.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.14'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +20.91%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0241'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.77'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +13.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.76'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0986'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.40%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.24'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.62%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '17.61'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("D49").Value = '1.511.67'
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.79'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.04%  '
